# Regenerate the "K" column (G) on the active sheet with freshly computed
# strikeout (K) values, replacing the old "Strike#" derived numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K (column G) value
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 3
    8  = 1
    9  = 2
    10 = 1
    12 = 2
    13 = 0
    14 = 2
    15 = 0
    16 = 4
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 3
    22 = 0
    23 = 2
    24 = 3
    25 = 0
    26 = 3
    27 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
